$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source formatting)
$ws.Range('D5,D6,D8,D10,D11,D15,D20,D21,D22,D23,D27,D32,D33,D34,D35,D37,D38,D40,D41,D42,D43,D44,D46,D47,D48,D49,D50,D51').NumberFormat = "@"

$ws.Range("D2").Value = '66.632.93'
$ws.Range("E2").Value = '  +3.76%  '
$ws.Range("D3").Value = '3.501.79'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '591.55'
$ws.Range("E5").Value = '  +3.21%  '
$ws.Range("D6").Value = '169.19'
$ws.Range("E6").Value = '  +2.90%  '
$ws.Range("D8").Value = '0.602'
$ws.Range("E8").Value = '  +8.63%  '
$ws.Range("D9").Value = '3.498.88'
$ws.Range("E9").Value = '  +1.76%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.128'
$ws.Range("E10").Value = '  +6.41%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = '7.34'
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("E12").Value = '  +3.46%  '
$ws.Range("D13").Value = '4.107.99'
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").Value = '28.33'
$ws.Range("E15").Value = '  +3.87%  '
$ws.Range("E16").Value = '  +2.65%  '
$ws.Range("D17").Value = '66.645.31'
$ws.Range("E17").Value = '  +3.69%  '
$ws.Range("D18").Value = '3.497.91'
$ws.Range("E18").Value = '  +1.71%  '
$ws.Range("E19").Value = '  +2.86%  '
$ws.Range("D20").Value = '14.17'
$ws.Range("E20").Value = '  +3.44%  '
$ws.Range("D21").Value = '392.88'
$ws.Range("E21").Value = '  +3.71%  '
$ws.Range("D22").Value = '7.98'
$ws.Range("E22").Value = '  +2.06%  '
$ws.Range("D23").Value = '73.20'
$ws.Range("E23").Value = '  +2.43%  '
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("E25").Value = '  +3.33%  '
$ws.Range("E26").Value = '  +5.08%  '
$ws.Range("D27").Value = '10.21'
$ws.Range("E27").Value = '  +6.64%  '
$ws.Range("E28").Value = '  +2.05%  '
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("E30").Value = '  +4.39%  '
$ws.Range("E31").Value = '  +4.67%  '
$ws.Range("D32").Value = '2.08'
$ws.Range("E32").Value = '  +3.07%  '
$ws.Range("D33").Value = '23.64'
$ws.Range("E33").Value = '  +2.73%  '
$ws.Range("D34").Value = '7.42'
$ws.Range("E34").Value = '  +3.97%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  +7.76%  '
$ws.Range("D37").Value = '162.42'
$ws.Range("E37").Value = '  +1.58%  '
$ws.Range("D38").Value = '0.885'
$ws.Range("E38").Value = '  +2.53%  '
$ws.Range("E39").Value = '  +4.47%  '
$ws.Range("D40").Value = '6.81'
$ws.Range("E40").Value = '  +5.06%  '
$ws.Range("D41").Value = '27.60'
$ws.Range("E41").Value = '  +4.31%  '
$ws.Range("D42").Value = '4.68'
$ws.Range("E42").Value = '  +5.64%  '
$ws.Range("D43").Value = '0.0746'
$ws.Range("E43").Value = '  +2.31%  '
$ws.Range("D44").Value = '26.48'
$ws.Range("E44").Value = '  +1.33%  '
$ws.Range("D45").Value = '2.797.02'
$ws.Range("E45").Value = '  -0.66%  '
$ws.Range("D46").Value = '43.25'
$ws.Range("E46").Value = '  +0.68%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '0.0312'
$ws.Range("E47").Value = '  +0.86%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").Value = '2.52'
$ws.Range("E48").Value = '  +1.14%  '
$ws.Range("D49").Value = '350.65'
$ws.Range("E49").Value = '  +4.73%  '
$ws.Range("D50").Value = '1.11'
$ws.Range("E50").Value = '  +4.98%  '
$ws.Range("D51").Value = '33.73'
$ws.Range("E51").Value = '  +12.03%  '
